# Rerun and summarise models without urban landuse: regenerate the
# per-model coefficient/p-value summary sheets, dropping the
# 'Education[T.Secondary+BAC]' row (models were rerun without it) and
# renaming each sheet to match the new run's summary-table id.
$wb = $excel.ActiveWorkbook

# Sheet 1: summ46266556 -> summ27923346
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ27923346"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 3150.295450997176
$ws.Cells.Item(2,3).Value = 0.009776264340039099
$ws.Cells.Item(3,2).Value = -18.24240633719774
$ws.Cells.Item(3,3).Value = 0.9583086591474073
$ws.Cells.Item(4,2).Value = 916.5534917073849
$ws.Cells.Item(4,3).Value = 0.008478460690204805
$ws.Cells.Item(5,2).Value = -644.5451410012984
$ws.Cells.Item(5,3).Value = 0.3048016915143009
$ws.Cells.Item(6,2).Value = -16.30874735232493
$ws.Cells.Item(6,3).Value = 0.6831481783201294
$ws.Cells.Item(7,2).Value = -1169.932530671134
$ws.Cells.Item(7,3).Value = 0.0000000000000000000000001174378199164749
$ws.Cells.Item(8,2).Value = -12.44934823552679
$ws.Cells.Item(8,3).Value = 0.01436589324313701
$ws.Cells.Item(9,2).Value = 196.3455869303829
$ws.Cells.Item(9,3).Value = 0.0000000006905878461149653
$ws.Cells.Item(10,2).Value = 458.717099332233
$ws.Cells.Item(10,3).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000011632692268253
$ws.Cells.Item(11,2).Value = -0.01597244745192818
$ws.Cells.Item(11,3).Value = 0.2005026065271671
$ws.Cells.Item(12,2).Value = 0.00001772450210254636
$ws.Cells.Item(12,3).Value = 0.3734494745773345
$ws.Cells.Item(13,2).Value = -3.59121198866437
$ws.Cells.Item(13,3).Value = 0.5412114861393367
$ws.Cells.Item(14,2).Value = 0.5160732097115259
$ws.Cells.Item(14,3).Value = 0.9468056921014898
$ws.Cells.Item(15,2).Value = -839.1051773427689
$ws.Cells.Item(15,3).Value = 0.2583387900435328
$ws.Cells.Item(16,2).Value = 1446.151746091651
$ws.Cells.Item(16,3).Value = 0.01007142560682553
$ws.Cells.Item(17,2).Value = 2600.682794629907
$ws.Cells.Item(17,3).Value = 0.00009194117799245378

# Sheet 2: summ46827709 -> summ28441840
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ28441840"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 3869.810839923689
$ws.Cells.Item(2,3).Value = 0.001673808211402191
$ws.Cells.Item(3,2).Value = -271.5889013116354
$ws.Cells.Item(3,3).Value = 0.4145267625205715
$ws.Cells.Item(4,2).Value = 692.2158204911979
$ws.Cells.Item(4,3).Value = 0.03671155568038143
$ws.Cells.Item(5,2).Value = -721.3772719740043
$ws.Cells.Item(5,3).Value = 0.2399555974623473
$ws.Cells.Item(6,2).Value = -0.8344529182729801
$ws.Cells.Item(6,3).Value = 0.9834672929790369
$ws.Cells.Item(7,2).Value = -1264.381844595368
$ws.Cells.Item(7,3).Value = 0.00000000000000000000000000001812640351384437
$ws.Cells.Item(8,2).Value = -10.25278296491948
$ws.Cells.Item(8,3).Value = 0.04358940833199066
$ws.Cells.Item(9,2).Value = 209.7314307391883
$ws.Cells.Item(9,3).Value = 0.00000000004751632160416919
$ws.Cells.Item(10,2).Value = 437.5434320866282
$ws.Cells.Item(10,3).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000006415061144275021
$ws.Cells.Item(11,2).Value = -0.01303405203394206
$ws.Cells.Item(11,3).Value = 0.295210593862347
$ws.Cells.Item(12,2).Value = -0.0000002304807811209506
$ws.Cells.Item(12,3).Value = 0.9906879965991139
$ws.Cells.Item(13,2).Value = -0.9285618225887653
$ws.Cells.Item(13,3).Value = 0.8754773220295486
$ws.Cells.Item(14,2).Value = -2.684507709693042
$ws.Cells.Item(14,3).Value = 0.7339175694895934
$ws.Cells.Item(15,2).Value = -189.7586738785022
$ws.Cells.Item(15,3).Value = 0.7996868647417941
$ws.Cells.Item(16,2).Value = 1018.4817947471
$ws.Cells.Item(16,3).Value = 0.07064517052578928
$ws.Cells.Item(17,2).Value = 2466.549664187867
$ws.Cells.Item(17,3).Value = 0.0002109704331034946

# Sheet 3: summ47480502 -> summ28991285
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ28991285"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 4073.956595744447
$ws.Cells.Item(2,3).Value = 0.0010053092258082
$ws.Cells.Item(3,2).Value = -378.8511596687486
$ws.Cells.Item(3,3).Value = 0.2640207204107247
$ws.Cells.Item(4,2).Value = 548.3152453555134
$ws.Cells.Item(4,3).Value = 0.104175914896794
$ws.Cells.Item(5,2).Value = -991.5613403986247
$ws.Cells.Item(5,3).Value = 0.12074448981108
$ws.Cells.Item(6,2).Value = -23.80719777661712
$ws.Cells.Item(6,3).Value = 0.5559700416268645
$ws.Cells.Item(7,2).Value = -1289.237796458542
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000003022767767492729
$ws.Cells.Item(8,2).Value = -15.08859319970821
$ws.Cells.Item(8,3).Value = 0.003057971560379366
$ws.Cells.Item(9,2).Value = 194.5845802500208
$ws.Cells.Item(9,3).Value = 0.0000000008618123223596253
$ws.Cells.Item(10,2).Value = 445.4080578700086
$ws.Cells.Item(10,3).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000005519965389587618
$ws.Cells.Item(11,2).Value = -0.01507732152329226
$ws.Cells.Item(11,3).Value = 0.2295105775058248
$ws.Cells.Item(12,2).Value = 0.00001232168260088399
$ws.Cells.Item(12,3).Value = 0.5350317186624872
$ws.Cells.Item(13,2).Value = 0.4148868275790065
$ws.Cells.Item(13,3).Value = 0.9449248441077067
$ws.Cells.Item(14,2).Value = 0.3345619567375508
$ws.Cells.Item(14,3).Value = 0.9663428979848452
$ws.Cells.Item(15,2).Value = -826.9488287714722
$ws.Cells.Item(15,3).Value = 0.2702935480210087
$ws.Cells.Item(16,2).Value = 884.040944328054
$ws.Cells.Item(16,3).Value = 0.1212387102891517
$ws.Cells.Item(17,2).Value = 1966.502439952703
$ws.Cells.Item(17,3).Value = 0.003282258534097854

# Sheet 4: summ48067171 -> summ29609118
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ29609118"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 4516.117914734119
$ws.Cells.Item(2,3).Value = 0.0002389290551782089
$ws.Cells.Item(3,2).Value = -476.2431245611814
$ws.Cells.Item(3,3).Value = 0.1685329902581101
$ws.Cells.Item(4,2).Value = 461.036927010548
$ws.Cells.Item(4,3).Value = 0.1805485571298459
$ws.Cells.Item(5,2).Value = -1212.66178595309
$ws.Cells.Item(5,3).Value = 0.05731521134202013
$ws.Cells.Item(6,2).Value = -1.755984626835428
$ws.Cells.Item(6,3).Value = 0.9649170453082656
$ws.Cells.Item(7,2).Value = -1283.052834203978
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000002551675662354845
$ws.Cells.Item(8,2).Value = -6.994658485573814
$ws.Cells.Item(8,3).Value = 0.1683051639736834
$ws.Cells.Item(9,2).Value = 208.6866692050559
$ws.Cells.Item(9,3).Value = 0.00000000004956754053271539
$ws.Cells.Item(10,2).Value = 439.3868608056421
$ws.Cells.Item(10,3).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000003525693593093653
$ws.Cells.Item(11,2).Value = -0.01640517546751384
$ws.Cells.Item(11,3).Value = 0.1829640246654735
$ws.Cells.Item(12,2).Value = 0.000003654049568006849
$ws.Cells.Item(12,3).Value = 0.8518162901474289
$ws.Cells.Item(13,2).Value = -6.180853271746052
$ws.Cells.Item(13,3).Value = 0.2951961161562838
$ws.Cells.Item(14,2).Value = -6.515099808102988
$ws.Cells.Item(14,3).Value = 0.4079667655581831
$ws.Cells.Item(15,2).Value = -394.3335211703228
$ws.Cells.Item(15,3).Value = 0.6000833803343171
$ws.Cells.Item(16,2).Value = 1506.878598690513
$ws.Cells.Item(16,3).Value = 0.007517999154391556
$ws.Cells.Item(17,2).Value = 2502.936711806895
$ws.Cells.Item(17,3).Value = 0.0001729850711549536

# Sheet 5: summ48642504 -> summ30221999
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ30221999"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 1446.299493754334
$ws.Cells.Item(2,3).Value = 0.2381443923207373
$ws.Cells.Item(3,2).Value = -34.06470520904108
$ws.Cells.Item(3,3).Value = 0.9201052271802505
$ws.Cells.Item(4,2).Value = 894.5292302456917
$ws.Cells.Item(4,3).Value = 0.008097296364483149
$ws.Cells.Item(5,2).Value = -426.2692723303691
$ws.Cells.Item(5,3).Value = 0.502606721336468
$ws.Cells.Item(6,2).Value = 34.19797266218511
$ws.Cells.Item(6,3).Value = 0.3946594511080889
$ws.Cells.Item(7,2).Value = -1284.384221631192
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000001837035500951122
$ws.Cells.Item(8,2).Value = -10.00327256496355
$ws.Cells.Item(8,3).Value = 0.0486489763246992
$ws.Cells.Item(9,2).Value = 192.9623232465693
$ws.Cells.Item(9,3).Value = 0.000000001726167498549023
$ws.Cells.Item(10,2).Value = 463.2408616894999
$ws.Cells.Item(10,3).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000001398170538947005
$ws.Cells.Item(11,2).Value = -0.007270585134981284
$ws.Cells.Item(11,3).Value = 0.5575473571180802
$ws.Cells.Item(12,2).Value = 0.000006903389413449208
$ws.Cells.Item(12,3).Value = 0.7269307712766242
$ws.Cells.Item(13,2).Value = 1.675756669636181
$ws.Cells.Item(13,3).Value = 0.7758922810611703
$ws.Cells.Item(14,2).Value = 15.11883428260446
$ws.Cells.Item(14,3).Value = 0.05300876866540097
$ws.Cells.Item(15,2).Value = -825.9386047616526
$ws.Cells.Item(15,3).Value = 0.2675981597057322
$ws.Cells.Item(16,2).Value = 835.1062502983252
$ws.Cells.Item(16,3).Value = 0.1377997916877778
$ws.Cells.Item(17,2).Value = 3009.104483808304
$ws.Cells.Item(17,3).Value = 0.00000547039104792275

# Sheet 6: summ49252633 -> summ30777322
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ30777322"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 3443.430644935967
$ws.Cells.Item(2,3).Value = 0.005069156154890366
$ws.Cells.Item(3,2).Value = -21.40089829865605
$ws.Cells.Item(3,3).Value = 0.9496163930689603
$ws.Cells.Item(4,2).Value = 940.772234337272
$ws.Cells.Item(4,3).Value = 0.005257900663798894
$ws.Cells.Item(5,2).Value = -300.9705144156941
$ws.Cells.Item(5,3).Value = 0.6314474167098094
$ws.Cells.Item(6,2).Value = 35.28362878316373
$ws.Cells.Item(6,3).Value = 0.3759260353353358
$ws.Cells.Item(7,2).Value = -1301.943809082521
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000000260325446372386
$ws.Cells.Item(8,2).Value = -12.3561129759599
$ws.Cells.Item(8,3).Value = 0.01537193935378312
$ws.Cells.Item(9,2).Value = 203.6896492126291
$ws.Cells.Item(9,3).Value = 0.0000000001417961629807153
$ws.Cells.Item(10,2).Value = 472.0404348105596
$ws.Cells.Item(10,3).Value = 0.00000000000000000000000000000000000000000000000000000000000000000000000000000000002050375944186655
$ws.Cells.Item(11,2).Value = -0.008495949676624547
$ws.Cells.Item(11,3).Value = 0.4927387019928288
$ws.Cells.Item(12,2).Value = 0.00001337456453052258
$ws.Cells.Item(12,3).Value = 0.4964744485318819
$ws.Cells.Item(13,2).Value = -1.880943647627465
$ws.Cells.Item(13,3).Value = 0.7496046869181153
$ws.Cells.Item(14,2).Value = -1.639813078790691
$ws.Cells.Item(14,3).Value = 0.8341224785670156
$ws.Cells.Item(15,2).Value = -1114.247908902239
$ws.Cells.Item(15,3).Value = 0.137323970253333
$ws.Cells.Item(16,2).Value = 706.8632189131893
$ws.Cells.Item(16,3).Value = 0.208510043748723
$ws.Cells.Item(17,2).Value = 2877.490746520743
$ws.Cells.Item(17,3).Value = 0.00001303318015107235

# Sheet 7: summ49809437 -> summ31378646
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ31378646"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 5040.152296276061
$ws.Cells.Item(2,3).Value = 0.00004613696924477603
$ws.Cells.Item(3,2).Value = -231.9149609138465
$ws.Cells.Item(3,3).Value = 0.4913497740692125
$ws.Cells.Item(4,2).Value = 652.8137086684517
$ws.Cells.Item(4,3).Value = 0.05139812980044114
$ws.Cells.Item(5,2).Value = -419.8506668429701
$ws.Cells.Item(5,3).Value = 0.5368722761532314
$ws.Cells.Item(6,2).Value = 37.01429314592936
$ws.Cells.Item(6,3).Value = 0.3570120603971458
$ws.Cells.Item(7,2).Value = -1254.174468699154
$ws.Cells.Item(7,3).Value = 0.00000000000000000000000000005348463021452516
$ws.Cells.Item(8,2).Value = -13.27596888812777
$ws.Cells.Item(8,3).Value = 0.009050818721296608
$ws.Cells.Item(9,2).Value = 179.4510919546243
$ws.Cells.Item(9,3).Value = 0.00000002045508007481905
$ws.Cells.Item(10,2).Value = 429.8877794223127
$ws.Cells.Item(10,3).Value = 0.000000000000000000000000000000000000000000000000000000000000000000009416178212570482
$ws.Cells.Item(11,2).Value = -0.0118383171999337
$ws.Cells.Item(11,3).Value = 0.3429693559075786
$ws.Cells.Item(12,2).Value = 0.000001428876289881439
$ws.Cells.Item(12,3).Value = 0.9423335963109374
$ws.Cells.Item(13,2).Value = -9.43532505182123
$ws.Cells.Item(13,3).Value = 0.1113556893557027
$ws.Cells.Item(14,2).Value = -8.100798808265033
$ws.Cells.Item(14,3).Value = 0.3036632842234305
$ws.Cells.Item(15,2).Value = -283.0857756402536
$ws.Cells.Item(15,3).Value = 0.7056797036997845
$ws.Cells.Item(16,2).Value = 1253.957235781088
$ws.Cells.Item(16,3).Value = 0.02652821106935778
$ws.Cells.Item(17,2).Value = 2453.164043624271
$ws.Cells.Item(17,3).Value = 0.0002321605955677006

# Sheet 8: summ50431824 -> summ31961491
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ31961491"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 3186.255906703884
$ws.Cells.Item(2,3).Value = 0.009376418239962681
$ws.Cells.Item(3,2).Value = -131.4208755142864
$ws.Cells.Item(3,3).Value = 0.6990667567553599
$ws.Cells.Item(4,2).Value = 787.309000925144
$ws.Cells.Item(4,3).Value = 0.02020873084015517
$ws.Cells.Item(5,2).Value = -672.3881389801609
$ws.Cells.Item(5,3).Value = 0.286853922904057
$ws.Cells.Item(6,2).Value = -8.40117022580845
$ws.Cells.Item(6,3).Value = 0.8343439754867134
$ws.Cells.Item(7,2).Value = -1293.769391790179
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000001555349400056726
$ws.Cells.Item(8,2).Value = -9.241522064155964
$ws.Cells.Item(8,3).Value = 0.07268338290456716
$ws.Cells.Item(9,2).Value = 219.1951469040575
$ws.Cells.Item(9,3).Value = 0.000000000007803094451512181
$ws.Cells.Item(10,2).Value = 440.4712192500074
$ws.Cells.Item(10,3).Value = 0.000000000000000000000000000000000000000000000000000000000000000000000009460761325384719
$ws.Cells.Item(11,2).Value = -0.01484302572613685
$ws.Cells.Item(11,3).Value = 0.2353998068173889
$ws.Cells.Item(12,2).Value = 0.000004576451543142907
$ws.Cells.Item(12,3).Value = 0.8195560472127228
$ws.Cells.Item(13,2).Value = -1.3257264074967
$ws.Cells.Item(13,3).Value = 0.8229963225996356
$ws.Cells.Item(14,2).Value = 1.039171924025017
$ws.Cells.Item(14,3).Value = 0.894747809007852
$ws.Cells.Item(15,2).Value = 30.31080216696252
$ws.Cells.Item(15,3).Value = 0.968061113843451
$ws.Cells.Item(16,2).Value = 1266.195986288083
$ws.Cells.Item(16,3).Value = 0.02539201207562448
$ws.Cells.Item(17,2).Value = 2562.380205185103
$ws.Cells.Item(17,3).Value = 0.0001234915155474198

# Sheet 9: summ51010828 -> summ32540268
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ32540268"
$ws.Rows.Item(4).Delete()  # drop Education[T.Secondary+BAC]
$ws.Cells.Item(2,2).Value = 2818.291643881035
$ws.Cells.Item(2,3).Value = 0.02215352707661785
$ws.Cells.Item(3,2).Value = -23.08605453740054
$ws.Cells.Item(3,3).Value = 0.9447766091007832
$ws.Cells.Item(4,2).Value = 867.6129484516017
$ws.Cells.Item(4,3).Value = 0.008973532910162497
$ws.Cells.Item(5,2).Value = -394.7926936912777
$ws.Cells.Item(5,3).Value = 0.5355413349735151
$ws.Cells.Item(6,2).Value = 8.18424456548299
$ws.Cells.Item(6,3).Value = 0.8376145781308905
$ws.Cells.Item(7,2).Value = -1345.981292894136
$ws.Cells.Item(7,3).Value = 0.000000000000000000000000000000004685385501332331
$ws.Cells.Item(8,2).Value = -12.87598824184023
$ws.Cells.Item(8,3).Value = 0.01187088983100395
$ws.Cells.Item(9,2).Value = 186.9715633104643
$ws.Cells.Item(9,3).Value = 0.000000005160677911742338
$ws.Cells.Item(10,2).Value = 452.6910713923112
$ws.Cells.Item(10,3).Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000006844303932637915
$ws.Cells.Item(11,2).Value = -0.01688001378353536
$ws.Cells.Item(11,3).Value = 0.1783741741667743
$ws.Cells.Item(12,2).Value = 0.000008302205966719945
$ws.Cells.Item(12,3).Value = 0.6766883550431747
$ws.Cells.Item(13,2).Value = 0.09723833018586148
$ws.Cells.Item(13,3).Value = 0.9869218186546763
$ws.Cells.Item(14,2).Value = 6.461841752610383
$ws.Cells.Item(14,3).Value = 0.4137053948814322
$ws.Cells.Item(15,2).Value = -1292.858116618074
$ws.Cells.Item(15,3).Value = 0.08739961888755853
$ws.Cells.Item(16,2).Value = 1409.777663765177
$ws.Cells.Item(16,3).Value = 0.01275658825929554
$ws.Cells.Item(17,2).Value = 2767.171573656149
$ws.Cells.Item(17,3).Value = 0.00003324845845332639

